# Applies the Sagittarius_Profits leve-profit value refresh across all sheets.
# Mirrors a scheduled data-refresh run: column H/I/J/K/L/M/N numeric values
# (currentAveragePrice*, LevePrice*, LeveProfit*) are updated per leve row;
# a few rows gain newly-populated profit cells, and a few rows lose cells
# that no longer carry a value (cleared, not zeroed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4356.9
$ws.Range("J17").Value = 4403.6895
$ws.Range("L17").Value = 13211.0685
$ws.Range("N17").Value = -13547.0685
$ws.Range("H40").Value = 2298.88
$ws.Range("J40").Value = 2509.3333
$ws.Range("L40").Value = 2509.3333
$ws.Range("N40").Value = -2859.3333
$ws.Range("H64").Value = 5500
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996
$ws.Range("H67").Value = 5500
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216
$ws.Range("H100").Value = 1735.7142
$ws.Range("I100").Value = 825
$ws.Range("J100").Value = 2100
$ws.Range("K100").Value = 825
$ws.Range("L100").Value = 2100
$ws.Range("M100").Value = -284
$ws.Range("N100").Value = -3182
$ws.Range("H138").Value = 3420.9194
$ws.Range("J138").Value = 3793.6223
$ws.Range("L138").Value = 11380.8669
$ws.Range("N138").Value = -21660.8669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2104.077
$ws.Range("I122").Value = 1377
$ws.Range("K122").Value = 4131
$ws.Range("M122").Value = -1681
$ws.Range("H132").Value = 1616.1724
$ws.Range("I132").Value = 1576.4814
$ws.Range("K132").Value = 4729.4442
$ws.Range("M132").Value = -2199.4442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19760
$ws.Range("H27").Value = 20000
$ws.Range("I27").Value = 20000
$ws.Range("K27").Value = 20000
$ws.Range("M27").Value = -19808
$ws.Range("H58").Value = 1866.3125
$ws.Range("I58").Value = 1872.75
$ws.Range("J58").Value = 1859.875
$ws.Range("K58").Value = 1872.75
$ws.Range("L58").Value = 1859.875
$ws.Range("M58").Value = -1669.75
$ws.Range("N58").Value = -2265.875
$ws.Range("H105").Value = 3443
$ws.Range("I105").Value = 2588.6428
$ws.Range("K105").Value = 2588.6428
$ws.Range("M105").Value = -841.6428000000001
$ws.Range("H132").Value = 1601.1282
$ws.Range("I132").Value = 1305.8387
$ws.Range("K132").Value = 3917.5161
$ws.Range("M132").Value = -1387.5161
$ws.Range("H134").Value = 1939.6818
$ws.Range("I134").Value = 1971.6666
$ws.Range("J134").Value = 1795.75
$ws.Range("K134").Value = 5914.9998
$ws.Range("L134").Value = 5387.25
$ws.Range("M134").Value = -3379.9998
$ws.Range("N134").Value = -10457.25
$ws.Range("H136").Value = 1866.3125
$ws.Range("I136").Value = 1872.75
$ws.Range("J136").Value = 1859.875
$ws.Range("K136").Value = 5618.25
$ws.Range("L136").Value = 5579.625
$ws.Range("M136").Value = -3068.25
$ws.Range("N136").Value = -10679.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8000139.5
$ws.Range("I4").Value = 7857292.5
$ws.Range("K4").Value = 23571877.5
$ws.Range("M4").Value = -23571765.5
$ws.Range("H34").Value = 500
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H38").Value = 38.333332
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 15
$ws.Range("K38").Value = 150
$ws.Range("L38").Value = 45
$ws.Range("M38").Value = 197
$ws.Range("N38").Value = -739
$ws.Range("H64").Value = 19596.666
$ws.Range("J64").Value = 19596.666
$ws.Range("L64").Value = 58789.99800000001
$ws.Range("N64").Value = -59329.99800000001
$ws.Range("H67").Value = 19596.666
$ws.Range("J67").Value = 19596.666
$ws.Range("L67").Value = 58789.99800000001
$ws.Range("N67").Value = -60661.99800000001
$ws.Range("H94").Value = 15806.125
$ws.Range("J94").Value = 19408.334
$ws.Range("L94").Value = 58225.00199999999
$ws.Range("N94").Value = -59577.00199999999
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H140").Value = 5996.6665
$ws.Range("I140").Value = 1495.25
$ws.Range("K140").Value = 4485.75
$ws.Range("M140").Value = 694.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H122").Value = 4605.16
$ws.Range("I122").Value = 3603.8125
$ws.Range("J122").Value = 6385.3335
$ws.Range("K122").Value = 10811.4375
$ws.Range("L122").Value = 19156.0005
$ws.Range("M122").Value = -8361.4375
$ws.Range("N122").Value = -24056.0005
$ws.Range("H132").Value = 1214.5
$ws.Range("I132").Value = 930
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 2790
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -260
$ws.Range("N132").Value = -9557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 142500
$ws.Range("I2").Value = 142500
$ws.Range("K2").Value = 142500
$ws.Range("M2").Value = -142388
$ws.Range("H22").Value = 1600
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 1800
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1800
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -2390
$ws.Range("H27").Value = 1600
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 1800
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 1800
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -2014
$ws.Range("H46").Value = 58238.555
$ws.Range("I46").Value = 250800
$ws.Range("J46").Value = 3221
$ws.Range("K46").Value = 250800
$ws.Range("L46").Value = 3221
$ws.Range("M46").Value = -250612
$ws.Range("N46").Value = -3597
$ws.Range("H55").Value = 522
$ws.Range("I55").Value = 437
$ws.Range("K55").Value = 437
$ws.Range("M55").Value = -264
$ws.Range("H122").Value = 8177.1816
$ws.Range("I122").Value = 8707.286
$ws.Range("J122").Value = 7249.5
$ws.Range("K122").Value = 26121.858
$ws.Range("L122").Value = 21748.5
$ws.Range("M122").Value = -23671.858
$ws.Range("N122").Value = -26648.5
$ws.Range("H136").Value = 3257.9546
$ws.Range("J136").Value = 4319.6
$ws.Range("L136").Value = 12958.8
$ws.Range("N136").Value = -18058.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 119900
$ws.Range("I57").Value = 119900
$ws.Range("K57").Value = 119900
$ws.Range("M57").Value = -119146
$ws.Range("H100").Value = 16670824
$ws.Range("I100").Value = 20002990
$ws.Range("K100").Value = 40005980
$ws.Range("M100").Value = -40005439

